$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three use-cases (rogue DHCP detection, Password Guessing, STP manipulation) are
# promoted to immediately follow the "ARP spoofing" row (row 2), pushing the
# remaining rows down. Rewrite rows 3-11 with their new content accordingly.

# Row 3: CyberSec.Booster: Detection of a rogue DHCP server
$ws.Range("A3").Value = "CyberSec.Booster: Detection of a rogue DHCP server"
$ws.Range("B3").Value = "T1557.003"
$ws.Range("C3").Value = "SW_DAI-4"
$ws.Range("D3").Value = "cbfe07a0-f345-11ee-95ad-619443a476e1"

# Row 4: CyberSec.Booster: Password Guessing
$ws.Range("A4").Value = "CyberSec.Booster: Password Guessing"
$ws.Range("B4").Value = "T1110.001"
$ws.Range("C4").Value = "SSH-5 SEC_LOGIN-4"
$ws.Range("D4").Value = "d36fef60-f347-11ee-95ad-619443a476e1"

# Row 5: CyberSec.Booster: Manipulation of the STP protocol
$ws.Range("A5").Value = "CyberSec.Booster: Manipulation of the STP protocol"
$ws.Range("B5").Value = "T1498.001"
$ws.Range("C5").Value = "SPANTREE"
$ws.Range("D5").Value = "267c1370-f346-11ee-95ad-619443a476e1"

# Row 6: CyberSec.Booster: MAC Flooding
$ws.Range("A6").Value = "CyberSec.Booster: MAC Flooding"
$ws.Range("B6").Value = "T1498.001"
$ws.Range("C6").Value = "MACNOTIFY-6 PORT_SECURITY-2"
$ws.Range("D6").Value = "0aab2540-f347-11ee-95ad-619443a476e1"

# Row 7: CyberSec.Booster: Scanning ip blocks
$ws.Range("A7").Value = "CyberSec.Booster: Scanning ip blocks"
$ws.Range("B7").Value = "T1595.001"
$ws.Range("D7").ClearFormats()
$ws.Range("D7").Value = "eb158cb0-e54b-11ee-aad9-f582020d7fab"

# Row 8: CyberSec.Booster: Vulnerability Scanning
$ws.Range("A8").Value = "CyberSec.Booster: Vulnerability Scanning"
$ws.Range("B8").Value = "T1595.002"
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = "d884bc00-f343-11ee-95ad-619443a476e1"

# Row 9: CyberSec.Booster: DNS data exfiltration
$ws.Range("A9").Value = "CyberSec.Booster: DNS data exfiltration"
$ws.Range("B9").Value = "T1071.004"
$ws.Range("D9").Value = "8adc2400-f345-11ee-95ad-619443a476e1"
$ws.Range("D9").NumberFormat = "0.00E+00"

# Row 10: CyberSec.Booster: Telnet Port Activity
$ws.Range("A10").Value = "CyberSec.Booster: Telnet Port Activity"
$ws.Range("B10").Value = "TA0011"
$ws.Range("C10").ClearContents()
$ws.Range("D10").Value = "9610b910-f347-11ee-95ad-619443a476e1"

# Row 11: CyberSec.Booster: Manipulation of the DTP protocol
$ws.Range("A11").Value = "CyberSec.Booster: Manipulation of the DTP protocol"
$ws.Range("B11").Value = "T1557"
$ws.Range("C11").Value = "SWITCHPORT-5"
$ws.Range("D11").Value = "0007b730-f346-11ee-95ad-619443a476e1"

# Update the active selection to match the printed console output range.
$ws.Range("A10:D11").Select()
